$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @('B','C','D','E','G','H','I','L','N','O')
$data = @{}

$data[2] = @(0.9938049312760313, 0.212234597619414, 0.07909651800753181, 0.09333381263081897, 0.4601385920929175, 0.5938273185414999, 0.5382333918353446, 0.2163124842965516, 1.048868541489476, 2.064056326840614)
$data[3] = @(0.894141305142341, 0.2000957193538682, 0.07167157949231751, 0.0939493106743825, 0.4555348137965325, 0.5957503478763329, 0.542907236514079, 0.2066827801411506, 1.054573738535794, 2.05802185600578)
$data[4] = @(0.832981842795391, 0.1925862729588204, 0.06714643977471724, 0.09437578624672405, 0.4530833153082909, 0.5972481569939418, 0.5461293598405277, 0.2008767839549819, 1.058495131001621, 2.055765051049235)
$data[5] = @(0.808069165303408, 0.1895121866545821, 0.06531091400866273, 0.09456179339879434, 0.4521784940125855, 0.5979382656465191, 0.5475309138285525, 0.198537715716796, 1.060198482022784, 2.055209243041674)
$data[6] = @(0.8039331048468057, 0.1890009013980745, 0.06500664022220803, 0.09459341763476559, 0.4520339325552527, 0.598057674070489, 0.5477689840421363, 0.1981509441087326, 1.060487689558535, 2.05513891746736)
$data[7] = @(0.8326458178249823, 0.1925448708929736, 0.06712165078725718, 0.09437824534012229, 0.45307073145365, 0.597257141164647, 0.5461479034266503, 0.2008451293265381, 1.058517676226742, 2.055756082446464)
$data[8] = @(0.9594347287775236, 0.2080608702481186, 0.07652938455244396, 0.09353596108869233, 0.4584732215399612, 0.5944245788426912, 0.5397717447146633, 0.2129700751298031, 1.050748945290501, 2.061674810565876)
$data[9] = @(1.208280204660753, 0.2380354711260111, 0.09524686587317888, 0.09226935554342042, 0.4720539524081317, 0.5913857070663227, 0.530068323986729, 0.2375907078301083, 1.038828195093068, 2.084793412247052)
$data[10] = @(1.391173789890615, 0.2597748365876384, 0.1091649898296936, 0.09157338792730307, 0.4838668208338675, 0.5906875889252348, 0.5246526250357277, 0.2561921195655827, 1.032082932818483, 2.108829628485665)
$data[11] = @(1.474379093451091, 0.2696017679894851, 0.1155334299847794, 0.09130768662364552, 0.4896424857620332, 0.5907034524477268, 0.5225623152718271, 0.264765466392376, 1.029450030488618, 2.121302795404887)
$data[12] = @(1.505886220085927, 0.2733138310179299, 0.1179503292382122, 0.09121438833770412, 0.4918875947263217, 0.5907574150970873, 0.5218245620291597, 0.2680279312569809, 1.028515538682313, 2.126247845381954)
$data[13] = @(1.49910066587546, 0.272514783523917, 0.1174295708384108, 0.09123415639859189, 0.4914014881680657, 0.590743660330034, 0.5219810556561058, 0.2673245945052827, 1.028714018675601, 2.125172971972319)
$data[14] = @(1.476971232335984, 0.2699073468685924, 0.1157321634183859, 0.0912998642850571, 0.4898260289234457, 0.5907069308477162, 0.5225005409227705, 0.2650335527923318, 1.029371896627495, 2.121705181344737)
$data[15] = @(1.463416142883773, 0.2683090146313987, 0.1146931433546143, 0.09134106509564965, 0.4888685723478545, 0.5906906784070856, 0.5228257507394645, 0.2636322943072997, 1.029783006365776, 2.119609948259125)
$data[16] = @(1.385736096961352, 0.2591313498184888, 0.1087495409935428, 0.09159177615564751, 0.4834974719523046, 0.5906932624972114, 0.5247967539982454, 0.2556340653140268, 1.032263754880077, 2.108045482585482)
$data[17] = @(1.338082189382931, 0.2534850160134852, 0.1051128046401004, 0.09175861400752083, 0.4803055578262416, 0.5907802550330246, 0.5261016013070545, 0.250755886823228, 1.033897095851607, 2.101345532180403)
$data[18] = @(1.310673602694294, 0.2502315317249497, 0.1030245350455345, 0.0918593655977098, 0.4785074735495698, 0.5908616754572336, 0.5268872464198893, 0.2479605860896044, 1.034877552242577, 2.097636727702849)
$data[19] = @(1.30139369774173, 0.249128957298808, 0.1023180817899743, 0.09189430125258902, 0.4779051628861737, 0.5908946331375944, 0.527159282588272, 0.2470159516110186, 1.035216563181883, 2.096405851730111)
$data[20] = @(1.343154971754416, 0.2540866861285451, 0.1054995808486012, 0.09174035801090419, 0.4806414268732055, 0.590767746286204, 0.5259590612465104, 0.2512740914541922, 1.033718981119875, 2.102043760869492)
$data[21] = @(1.483471222521644, 0.2706734650839451, 0.1162305891467668, 0.09128036572015219, 0.4902872041335939, 0.5907164176327342, 0.5223464943992582, 0.265706055931318, 1.029176965626831, 2.122717734349294)
$data[22] = @(1.575169834979931, 0.2814602647082438, 0.1232748624457543, 0.09102238246162386, 0.4969293755933393, 0.5909623958633148, 0.5202991254840228, 0.2752309501020278, 1.026572931974457, 2.137521994166804)
$data[23] = @(1.526229715572413, 0.2757081213226229, 0.119512375011368, 0.09115617133637954, 0.4933533243964519, 0.5908055338217508, 0.5213631064186472, 0.270138881046222, 1.027929439871414, 2.129502265796248)
$data[24] = @(1.340861603214705, 0.2538146938926786, 0.1053247114170972, 0.09174859648867084, 0.4804894652715745, 0.5907733036582243, 0.5260233931262306, 0.2510397823866555, 1.033799377738717, 2.101727646005173)
$data[25] = @(1.140944645379761, 0.2299757161398475, 0.09015424437164654, 0.09257079424106607, 0.4680587818929496, 0.5919383830791816, 0.5323929203965143, 0.2308400281045664, 1.041699082759109, 2.077303362587998)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
